$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.2520848399290218
    "C2" = 0.7424843809194045
    "D2" = 1.903216275652063
    "E2" = 1.379571047699995
    "F2" = 1.369577110019525
    "B3" = 0.6974398184037287
    "C3" = 1.552624331747446
    "D3" = 7.423039619458927
    "E3" = 2.72452557695077
    "F3" = 2.659952944235514
    "B4" = 1.179297011432999
    "C4" = 1.950015878589636
    "D4" = 10.42997277502151
    "E4" = 3.22954683740947
    "F4" = 3.037055390127207
    "B5" = 0.9186370495837128
    "C5" = 1.942955877847222
    "D5" = 11.15520434258843
    "E5" = 3.339940769323377
    "F5" = 3.244399679230311
    "B6" = 1.01406298738989
    "C6" = 1.964591607092075
    "D6" = 11.42823783535527
    "E6" = 3.380567679451969
    "F6" = 3.259016548116396
    "B7" = 0.8734069636221409
    "C7" = 1.740984057691856
    "D7" = 10.78249254741741
    "E7" = 3.283670590576559
    "F7" = 3.206762837755778
    "B8" = 1.003986551470107
    "C8" = 1.777346341468097
    "D8" = 11.03323908461691
    "E8" = 3.321631991147862
    "F8" = 3.208769669857039
    "B9" = 0.8645380383347879
    "C9" = 1.931740234399347
    "D9" = 17.21370736197478
    "E9" = 4.148940510777996
    "F9" = 4.158075926358173
    "B10" = 0.03960286749493035
    "C10" = 1.355899730545265
    "D10" = 4.943417352971498
    "E10" = 2.223379714077534
    "F10" = 2.30694412459439
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
